# Commit: Natmi following Dr Hou advice
#
# The Il10 -> Il10rb LR-pair sheet is recomputed to include "ECs" as an
# additional Sending-cluster group (the grid becomes 4 sending clusters x
# 4 target clusters: ECs, FAPs, M2, sCs). Existing rows 2-13 get shifted
# sending-cluster labels/values and rows 14-17 are newly added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending cluster = ECs, Target cluster = ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il10"
$ws.Range("C2").Value = "Il10rb"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 36.26745966666667
$ws.Range("H2").Value = 108.802379
$ws.Range("I2").Value = 0.643144256662053
$ws.Range("J2").Value = 0.6431442566620529
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 86.19469199999999
$ws.Range("N2").Value = 258.584076
$ws.Range("O2").Value = 0.454900350325626
$ws.Range("P2").Value = 0.454900350325626
$ws.Range("Q2").Value = 3126.062515590756
$ws.Range("R2").Value = 28134.5626403168
$ws.Range("S2").Value = 0.2925665476654822
$ws.Range("T2").Value = 0.2925665476654822

# Row 3: Sending cluster = ECs, Target cluster = FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Il10"
$ws.Range("C3").Value = "Il10rb"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 36.26745966666667
$ws.Range("H3").Value = 108.802379
$ws.Range("I3").Value = 0.643144256662053
$ws.Range("J3").Value = 0.6431442566620529
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 9.289272666666667
$ws.Range("N3").Value = 27.867818
$ws.Range("O3").Value = 0.04902498393215361
$ws.Range("P3").Value = 0.04902498393215361
$ws.Range("Q3").Value = 336.8983217710025
$ws.Range("R3").Value = 3032.084895939022
$ws.Range("S3").Value = 0.03153013684891403
$ws.Range("T3").Value = 0.03153013684891402

# Row 4: Sending cluster = ECs, Target cluster = M2
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Il10"
$ws.Range("C4").Value = "Il10rb"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 36.26745966666667
$ws.Range("H4").Value = 108.802379
$ws.Range("I4").Value = 0.643144256662053
$ws.Range("J4").Value = 0.6431442566620529
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 88.33691399999999
$ws.Range("N4").Value = 265.010742
$ws.Range("O4").Value = 0.4662061223594221
$ws.Range("P4").Value = 0.4662061223594221
$ws.Range("Q4").Value = 3203.755465572802
$ws.Range("R4").Value = 28833.79919015522
$ws.Range("S4").Value = 0.2998377900161486
$ws.Range("T4").Value = 0.2998377900161486

# Row 5: Sending cluster = ECs, Target cluster = sCs
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Il10"
$ws.Range("C5").Value = "Il10rb"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 36.26745966666667
$ws.Range("H5").Value = 108.802379
$ws.Range("I5").Value = 0.643144256662053
$ws.Range("J5").Value = 0.6431442566620529
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.659503
$ws.Range("N5").Value = 16.978509
$ws.Range("O5").Value = 0.02986854338279823
$ws.Range("P5").Value = 0.02986854338279823
$ws.Range("Q5").Value = 205.255796785879
$ws.Range("R5").Value = 1847.302171072911
$ws.Range("S5").Value = 0.01920978213150805
$ws.Range("T5").Value = 0.01920978213150805

# Row 6: Sending cluster = FAPs, Target cluster = ECs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Il10"
$ws.Range("C6").Value = "Il10rb"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.209141666666667
$ws.Range("H6").Value = 9.627425000000001
$ws.Range("I6").Value = 0.05690889438359308
$ws.Range("J6").Value = 0.05690889438359308
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 86.19469199999999
$ws.Range("N6").Value = 258.584076
$ws.Range("O6").Value = 0.454900350325626
$ws.Range("P6").Value = 0.454900350325626
$ws.Range("Q6").Value = 276.6109775427
$ws.Range("R6").Value = 2489.4987978843
$ws.Range("S6").Value = 0.02588787599174055
$ws.Range("T6").Value = 0.02588787599174054

# Row 7: Sending cluster = FAPs, Target cluster = FAPs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Il10"
$ws.Range("C7").Value = "Il10rb"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.209141666666667
$ws.Range("H7").Value = 9.627425000000001
$ws.Range("I7").Value = 0.05690889438359308
$ws.Range("J7").Value = 0.05690889438359308
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 9.289272666666667
$ws.Range("N7").Value = 27.867818
$ws.Range("O7").Value = 0.04902498393215361
$ws.Range("P7").Value = 0.04902498393215361
$ws.Range("Q7").Value = 29.81059196762778
$ws.Range("R7").Value = 268.29532770865
$ws.Range("S7").Value = 0.002789957632752278
$ws.Range("T7").Value = 0.002789957632752277

# Row 8: Sending cluster = FAPs, Target cluster = M2
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Il10"
$ws.Range("C8").Value = "Il10rb"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.209141666666667
$ws.Range("H8").Value = 9.627425000000001
$ws.Range("I8").Value = 0.05690889438359308
$ws.Range("J8").Value = 0.05690889438359308
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 88.33691399999999
$ws.Range("N8").Value = 265.010742
$ws.Range("O8").Value = 0.4662061223594221
$ws.Range("P8").Value = 0.4662061223594221
$ws.Range("Q8").Value = 283.48567142215
$ws.Range("R8").Value = 2551.37104279935
$ws.Range("S8").Value = 0.02653127497833683
$ws.Range("T8").Value = 0.02653127497833682

# Row 9: Sending cluster = FAPs, Target cluster = sCs
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Il10"
$ws.Range("C9").Value = "Il10rb"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.209141666666667
$ws.Range("H9").Value = 9.627425000000001
$ws.Range("I9").Value = 0.05690889438359308
$ws.Range("J9").Value = 0.05690889438359308
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 5.659503
$ws.Range("N9").Value = 16.978509
$ws.Range("O9").Value = 0.02986854338279823
$ws.Range("P9").Value = 0.02986854338279823
$ws.Range("Q9").Value = 18.162146889925
$ws.Range("R9").Value = 163.459322009325
$ws.Range("S9").Value = 0.001699785780763433
$ws.Range("T9").Value = 0.001699785780763432

# Row 10: Sending cluster = M2, Target cluster = ECs
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Il10"
$ws.Range("C10").Value = "Il10rb"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 15.98639866666667
$ws.Range("H10").Value = 47.95919600000001
$ws.Range("I10").Value = 0.2834927116945642
$ws.Range("J10").Value = 0.2834927116945642
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 86.19469199999999
$ws.Range("N10").Value = 258.584076
$ws.Range("O10").Value = 0.454900350325626
$ws.Range("P10").Value = 0.454900350325626
$ws.Range("Q10").Value = 1377.942709262544
$ws.Range("R10").Value = 12401.4843833629
$ws.Range("S10").Value = 0.128960933864619
$ws.Range("T10").Value = 0.128960933864619

# Row 11: Sending cluster = M2, Target cluster = FAPs
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Il10"
$ws.Range("C11").Value = "Il10rb"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 15.98639866666667
$ws.Range("H11").Value = 47.95919600000001
$ws.Range("I11").Value = 0.2834927116945642
$ws.Range("J11").Value = 0.2834927116945642
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 9.289272666666667
$ws.Range("N11").Value = 27.867818
$ws.Range("O11").Value = 0.04902498393215361
$ws.Range("P11").Value = 0.04902498393215361
$ws.Range("Q11").Value = 148.5020161727031
$ws.Range("R11").Value = 1336.518145554328
$ws.Range("S11").Value = 0.01389822563570867
$ws.Range("T11").Value = 0.01389822563570866

# Row 12: Sending cluster = M2, Target cluster = M2
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Il10"
$ws.Range("C12").Value = "Il10rb"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 15.98639866666667
$ws.Range("H12").Value = 47.95919600000001
$ws.Range("I12").Value = 0.2834927116945642
$ws.Range("J12").Value = 0.2834927116945642
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 88.33691399999999
$ws.Range("N12").Value = 265.010742
$ws.Range("O12").Value = 0.4662061223594221
$ws.Range("P12").Value = 0.4662061223594221
$ws.Range("Q12").Value = 1412.189124187048
$ws.Range("R12").Value = 12709.70211768343
$ws.Range("S12").Value = 0.1321660378362804
$ws.Range("T12").Value = 0.1321660378362804

# Row 13: Sending cluster = M2, Target cluster = sCs
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Il10"
$ws.Range("C13").Value = "Il10rb"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 15.98639866666667
$ws.Range("H13").Value = 47.95919600000001
$ws.Range("I13").Value = 0.2834927116945642
$ws.Range("J13").Value = 0.2834927116945642
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 5.659503
$ws.Range("N13").Value = 16.978509
$ws.Range("O13").Value = 0.02986854338279823
$ws.Range("P13").Value = 0.02986854338279823
$ws.Range("Q13").Value = 90.475071213196
$ws.Range("R13").Value = 814.275640918764
$ws.Range("S13").Value = 0.008467514357956202
$ws.Range("T13").Value = 0.008467514357956202

# Row 14: Sending cluster = sCs, Target cluster = ECs
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Il10"
$ws.Range("C14").Value = "Il10rb"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.927863
$ws.Range("H14").Value = 2.783589
$ws.Range("I14").Value = 0.01645413725978976
$ws.Range("J14").Value = 0.01645413725978976
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 86.19469199999999
$ws.Range("N14").Value = 258.584076
$ws.Range("O14").Value = 0.454900350325626
$ws.Range("P14").Value = 0.454900350325626
$ws.Range("Q14").Value = 79.976865503196
$ws.Range("R14").Value = 719.791789528764
$ws.Range("S14").Value = 0.0074849928037843
$ws.Range("T14").Value = 0.0074849928037843

# Row 15: Sending cluster = sCs, Target cluster = FAPs
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Il10"
$ws.Range("C15").Value = "Il10rb"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.927863
$ws.Range("H15").Value = 2.783589
$ws.Range("I15").Value = 0.01645413725978976
$ws.Range("J15").Value = 0.01645413725978976
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 9.289272666666667
$ws.Range("N15").Value = 27.867818
$ws.Range("O15").Value = 0.04902498393215361
$ws.Range("P15").Value = 0.04902498393215361
$ws.Range("Q15").Value = 8.619172404311334
$ws.Range("R15").Value = 77.572551638802
$ws.Range("S15").Value = 0.0008066638147786432
$ws.Range("T15").Value = 0.0008066638147786431

# Row 16: Sending cluster = sCs, Target cluster = M2
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Il10"
$ws.Range("C16").Value = "Il10rb"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.927863
$ws.Range("H16").Value = 2.783589
$ws.Range("I16").Value = 0.01645413725978976
$ws.Range("J16").Value = 0.01645413725978976
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 88.33691399999999
$ws.Range("N16").Value = 265.010742
$ws.Range("O16").Value = 0.4662061223594221
$ws.Range("P16").Value = 0.4662061223594221
$ws.Range("Q16").Value = 81.964554034782
$ws.Range("R16").Value = 737.680986313038
$ws.Range("S16").Value = 0.007671019528656273
$ws.Range("T16").Value = 0.007671019528656273

# Row 17: Sending cluster = sCs, Target cluster = sCs
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Il10"
$ws.Range("C17").Value = "Il10rb"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.927863
$ws.Range("H17").Value = 2.783589
$ws.Range("I17").Value = 0.01645413725978976
$ws.Range("J17").Value = 0.01645413725978976
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 5.659503
$ws.Range("N17").Value = 16.978509
$ws.Range("O17").Value = 0.02986854338279823
$ws.Range("P17").Value = 0.02986854338279823
$ws.Range("Q17").Value = 5.251243432089
$ws.Range("R17").Value = 47.261190888801
$ws.Range("S17").Value = 0.0004914611125705474
$ws.Range("T17").Value = 0.0004914611125705473
